$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 'Clients'
$ws.Range("C3").Value = "Clients"
# 'Nom'
$ws.Range("C5").Value = "Nom"
$ws.Range("G5").Value = "Nom"
$ws.Range("I5").Value = "Nom"
$ws.Range("K5").Value = "Nom"
$ws.Range("M5").Value = "Nom"
# 'Cellulaire'
$ws.Range("C6").Value = "Cellulaire"
# 'Addresse'
$ws.Range("C7").Value = "Addresse"
# 'Code Postal'
$ws.Range("C11").Value = "Code Postal"
# 'Province'
$ws.Range("C9").Value = "Province"
# 'Pays'
$ws.Range("C10").Value = "Pays"
# 'Ville'
$ws.Range("C8").Value = "Ville"
# 'Commande'
$ws.Range("E3").Value = "Commande"
$ws.Range("O6").Value = "Commande"
# 'Croute'
$ws.Range("G3").Value = "Croute"
$ws.Range("E5").Value = "Croute"
# 'Sauce'
$ws.Range("K3").Value = "Sauce"
$ws.Range("E6").Value = "Sauce"
# 'Garniture'
$ws.Range("I3").Value = "Garniture"
$ws.Range("E7").Value = "Garniture"
# 'Prix'
$ws.Range("G6").Value = "Prix"
$ws.Range("I6").Value = "Prix"
$ws.Range("K6").Value = "Prix"
$ws.Range("M6").Value = "Prix"
# 'Commande Attente'
$ws.Range("O3").Value = "Commande Attente"
# 'Client'
$ws.Range("O5").Value = "Client"
$ws.Range("E9").Value = "Client"
# 'Taille'
$ws.Range("M3").Value = "Taille"
$ws.Range("E8").Value = "Taille"

$ws.Range("N18").Select()
